$wb = $excel.ActiveWorkbook

# Add the new "salesInvoice" worksheet after the last existing sheet (vendor)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "salesInvoice"

# Header row
$ws.Range("A1").Value = "sampleId"
$ws.Range("B1").Value = "userId"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "orgCode"
$ws.Range("E1").Value = "customerName"
$ws.Range("F1").Value = "store"
$ws.Range("G1").Value = "paymentMode"
$ws.Range("H1").Value = "productName"
$ws.Range("I1").Value = "quantity"
$ws.Range("J1").Value = "unitPrice"

# Row 2
$ws.Range("A2").Value = "Inv-1"
$ws.Range("B2").Value = "demo"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Pass@1234", "", "", "Pass@1234")
$ws.Range("D2").Value = "sid"
$ws.Range("E2").Value = "edward"
$ws.Range("F2").Value = "tirupur"
$ws.Range("G2").Value = "net banking"
$ws.Range("H2").Value = "zas Black Cardamom"
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 245

# Row 3
$ws.Range("A3").Value = "Inv-1"
$ws.Range("B3").Value = "demo"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Pass@1234", "", "", "Pass@1234")
$ws.Range("D3").Value = "sid"
$ws.Range("E3").Value = "edward"
$ws.Range("F3").Value = "tirupur"
$ws.Range("G3").Value = "net banking"
$ws.Range("H3").Value = "zas BlackGram"
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 124

# Row 4
$ws.Range("A4").Value = "Inv-2"
$ws.Range("B4").Value = "demo"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Pass@1234", "", "", "Pass@1234")
$ws.Range("D4").Value = "sid"
$ws.Range("E4").Value = "Zachary"
$ws.Range("F4").Value = "gowtham"
$ws.Range("G4").Value = "sbi"
$ws.Range("H4").Value = "zas DryChilli"
$ws.Range("I4").Value = 25
$ws.Range("J4").Value = 554

# Selection for the new sheet matches the diff (B3:G3 active range with A1 anchor in other sheets)
$ws.Range("B3:G3").Select()

Write-Host "Sheet created with" $wb.Worksheets.Count "total sheets"
